$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# Insert a new row at 34 (shifts everything below down by one) and populate it.
$ws.Rows(34).Insert()
$ws.Range("A34").Value = "axis"
$ws.Range("B34").Value = "Axis"

# New rows 55 and 56 (previously-unused row numbers) for the Newton's Second
# Law dialog entries.
$ws.Range("A55").Value = "newton_second_law_title"
$ws.Range("A56").Value = "newton_second_law_desc"
$ws.Range("B55").Value = "Newton's Second Law"
$ws.Range("C55").Value = 2
$ws.Range("B56").Value = $ws.Range("B20").Value2
$ws.Range("B55").VerticalAlignment = -4108

# Restore the selection/scroll state recorded in the edit.
$ws.Range("B52").Select() | Out-Null
